# Auto-generated edit script: updates crypto price/volume cells to match target snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.654.32"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.511.06"
$ws.Range("E3").Value = "  -2.44%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'587.30"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "'183.33"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("D7").Value = "3.499.98"
$ws.Range("E7").Value = "  -2.63%  "
$ws.Range("D8").Value = "'0.612"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.198"
$ws.Range("E10").Value = "  +6.92%  "
$ws.Range("D11").Value = "'0.644"
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").Value = "'53.93"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").Value = "'0.0000305"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "4.066.70"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").Value = "'19.28"
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").Value = "69.588.32"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "3.499.03"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "'531.75"
$ws.Range("E21").Value = "  +7.96%  "
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("D23").Value = "'18.33"
$ws.Range("E23").Value = "  -5.88%  "
$ws.Range("E24").Value = "  +5.11%  "
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").Value = "'95.27"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").Value = "'11.15"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "'9.08"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").Value = "'32.15"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "'7.28"
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "'63.83"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("D35").Value = "'544.71"
$ws.Range("E35").Value = "  -6.03%  "
$ws.Range("D36").Value = "'3.12"
$ws.Range("E36").Value = "  +6.38%  "
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "'38.05"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  -7.03%  "
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").Value = "3.353.27"
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("E44").Value = "  -4.55%  "
$ws.Range("D45").Value = "'2.97"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").Value = "'3.49"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").Value = "'0.0438"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").Value = "'8.97"
$ws.Range("E49").Value = "  -7.29%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "'138.31"
$ws.Range("E51").Value = "  +3.09%  "
